$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between rows 5 and 6 ---
$row5 = $ws.Range("F5:V5").Value2
$row6 = $ws.Range("F6:V6").Value2
$ws.Range("F5:V5").Value2 = $row6
$ws.Range("F6:V6").Value2 = $row5

# --- Swap the match data (columns F:V) between rows 58 and 59 ---
$row58 = $ws.Range("F58:V58").Value2
$row59 = $ws.Range("F59:V59").Value2
$ws.Range("F58:V58").Value2 = $row59
$ws.Range("F59:V59").Value2 = $row58

# --- Append two new match rows (62 and 63) ---

# Copy the formatting (styles) of the last existing row down into the two
# new rows so they look consistent with the rest of the table.
$ws.Range("A61").Copy()
$ws.Range("A62:A63").PasteSpecial(-4122)

$ws.Range("E61").Copy()
$ws.Range("E62:E63").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 62: Khenchela 2-0 El Bayadh
$ws.Cells.Item(62, 1).Value2 = 61
$ws.Cells.Item(62, 2).Value2 = "algeria"
$ws.Cells.Item(62, 3).Value2 = "ligue-1"
$ws.Cells.Item(62, 4).Value2 = "2023-2024"
$ws.Cells.Item(62, 5).Value2 = 45268.63541666666
$ws.Cells.Item(62, 6).Value2 = "Khenchela"
$ws.Cells.Item(62, 7).Value2 = 2
$ws.Cells.Item(62, 8).Value2 = "El Bayadh"
$ws.Cells.Item(62, 9).Value2 = 0
$ws.Cells.Item(62, 10).Value2 = 1.95
$ws.Cells.Item(62, 11).Value2 = "07/12/2023 07:08"
$ws.Cells.Item(62, 12).Value2 = 1.76
$ws.Cells.Item(62, 13).Value2 = "08/12/2023 15:10"
$ws.Cells.Item(62, 14).Value2 = 2.72
$ws.Cells.Item(62, 15).Value2 = "07/12/2023 07:08"
$ws.Cells.Item(62, 16).Value2 = 3.09
$ws.Cells.Item(62, 17).Value2 = "08/12/2023 15:10"
$ws.Cells.Item(62, 18).Value2 = 3.75
$ws.Cells.Item(62, 19).Value2 = "07/12/2023 07:08"
$ws.Cells.Item(62, 20).Value2 = 6.05
$ws.Cells.Item(62, 21).Value2 = "08/12/2023 15:06"
$ws.Cells.Item(62, 22).Value2 = "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-el-bayadh/vwCccbzs/"

# Row 63: Kabylie 1-0 Ben Aknoun
$ws.Cells.Item(63, 1).Value2 = 62
$ws.Cells.Item(63, 2).Value2 = "algeria"
$ws.Cells.Item(63, 3).Value2 = "ligue-1"
$ws.Cells.Item(63, 4).Value2 = "2023-2024"
$ws.Cells.Item(63, 5).Value2 = 45268.75
$ws.Cells.Item(63, 6).Value2 = "Kabylie"
$ws.Cells.Item(63, 7).Value2 = 1
$ws.Cells.Item(63, 8).Value2 = "Ben Aknoun"
$ws.Cells.Item(63, 9).Value2 = 0
$ws.Cells.Item(63, 10).Value2 = 1.32
$ws.Cells.Item(63, 11).Value2 = "07/12/2023 06:11"
$ws.Cells.Item(63, 12).Value2 = 1.29
$ws.Cells.Item(63, 13).Value2 = "08/12/2023 17:18"
$ws.Cells.Item(63, 14).Value2 = 3.7
$ws.Cells.Item(63, 15).Value2 = "07/12/2023 06:11"
$ws.Cells.Item(63, 16).Value2 = 4.8
$ws.Cells.Item(63, 17).Value2 = "08/12/2023 17:18"
$ws.Cells.Item(63, 18).Value2 = 6.88
$ws.Cells.Item(63, 19).Value2 = "07/12/2023 06:11"
$ws.Cells.Item(63, 20).Value2 = 13.66
$ws.Cells.Item(63, 21).Value2 = "08/12/2023 17:18"
$ws.Cells.Item(63, 22).Value2 = "https://www.betexplorer.com/football/algeria/ligue-1/kabylie-es-ben-aknoun/n92IhzkC/"
